# Applies: rename sheet Sheet1 -> LibraryAPI, update header row B1:E1 with
# new column names, and move the active selection to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "LibraryAPI"

$ws.Range("B1").Value = "BookName"
$ws.Range("C1").Value = "ISBN"
$ws.Range("D1").Value = "ISLE"
$ws.Range("E1").Value = "Author"

$ws.Range("F6").Select()
